# Convert the "Day" column from day-name text (Tuesday/Wednesday/Thursday/Friday)
# into numeric category codes (1/2/3/4), matching the move from named days to
# generic, user-editable categories.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tuesday (rows 2-11) -> category 1
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Wednesday (rows 12-17) -> category 2
for ($r = 12; $r -le 17; $r++) {
    $ws.Cells.Item($r, 2).Value = 2
}

# Thursday (rows 18-23) -> category 3
for ($r = 18; $r -le 23; $r++) {
    $ws.Cells.Item($r, 2).Value = 3
}

# Friday (rows 24-33) -> category 4
for ($r = 24; $r -le 33; $r++) {
    $ws.Cells.Item($r, 2).Value = 4
}

# Add a brand-new category (0), a schedule that begins a bit later than the
# others (08:40/09:25 for period 1) but otherwise mirrors the standard
# period/brunch/lunch layout.
$newRows = @(
    @(1, 0, "08:40", "09:25"),
    @(2, 0, "09:30", "10:15"),
    @(3, 0, "10:20", "11:10"),
    @("Brunch", 0, "11:10", "11:25"),
    @(4, 0, "11:30", "12:15"),
    @(5, 0, "12:20", "13:05"),
    @("Lunch", 0, "13:05", "13:45"),
    @(6, 0, "13:50", "14:35"),
    @(7, 0, "14:40", "15:25")
)

$row = 34
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.Value = $entry[2]
    $cellC.NumberFormat = "h:mm AM/PM"

    $cellD = $ws.Cells.Item($row, 4)
    $cellD.Value = $entry[3]
    $cellD.NumberFormat = "h:mm AM/PM"

    $row++
}

$null = $ws.Range("B43").Select()
